$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old column H (median_runtime duplicate column no longer needed after
# collapsing "praklasifikasi"/"klasifikasi" into a single "models" column)
$ws.Columns.Item(8).Delete()

# Header row
$ws.Range("A1").Value = "models"
$ws.Range("B1").Value = "accuracy"
$ws.Range("C1").Value = "recall"
$ws.Range("D1").Value = "precision"
$ws.Range("E1").Value = "f1_score"
$ws.Range("F1").Value = "average_runtime"
$ws.Range("G1").Value = "median_runtime"

# Row 2 - mamdani + mlp
$ws.Range("A2").Value = "['mamdani', 'mlp']"
$ws.Range("B2").Value = 0.806577480490524
$ws.Range("C2").Value = 0.806577480490524
$ws.Range("D2").Value = 0.8125072957494913
$ws.Range("E2").Value = 0.8069127034461214
$ws.Range("F2").Value = 0.2663744926452637
$ws.Range("G2").Value = 0.1319910287857056

# Row 3 - mamdani + lgbm
$ws.Range("A3").Value = "['mamdani', 'lgbm']"
$ws.Range("B3").Value = 0.8344481605351171
$ws.Range("C3").Value = 0.8344481605351171
$ws.Range("D3").Value = 0.8361176842296993
$ws.Range("E3").Value = 0.8340194950146642
$ws.Range("F3").Value = 0.3190764278173447
$ws.Range("G3").Value = 0.1499724388122559

# Row 4 - lr + mlp
$ws.Range("A4").Value = "['lr', 'mlp']"
$ws.Range("B4").Value = 0.7987736900780379
$ws.Range("C4").Value = 0.7987736900780379
$ws.Range("D4").Value = 0.8020056261781602
$ws.Range("E4").Value = 0.798970837990839
$ws.Range("F4").Value = 0.2872806489467621
$ws.Range("G4").Value = 0.1328588724136353

# Row 5 - lr + lgbm
$ws.Range("A5").Value = "['lr', 'lgbm']"
$ws.Range("B5").Value = 0.8255295429208472
$ws.Range("C5").Value = 0.8255295429208472
$ws.Range("D5").Value = 0.8261889021608652
$ws.Range("E5").Value = 0.8253038499871026
$ws.Range("F5").Value = 0.345859083533287
$ws.Range("G5").Value = 0.1665419340133667

# Row 6 - dt + mlp
$ws.Range("A6").Value = "['dt', 'mlp']"
$ws.Range("B6").Value = 0.8132664437012264
$ws.Range("C6").Value = 0.8132664437012264
$ws.Range("D6").Value = 0.8318556751371609
$ws.Range("E6").Value = 0.8127850574242346
$ws.Range("F6").Value = 0.2706254035234451
$ws.Range("G6").Value = 0.1317030191421509

# Row 7 - dt + lgbm
$ws.Range("A7").Value = "['dt', 'lgbm']"
$ws.Range("B7").Value = 0.8355629877369007
$ws.Range("C7").Value = 0.8355629877369007
$ws.Range("D7").Value = 0.849953478282931
$ws.Range("E7").Value = 0.8348438778681438
$ws.Range("F7").Value = 0.3602471113204956
$ws.Range("G7").Value = 0.1784284114837646

# Row 8 - lgbm (single model row)
$ws.Range("A8").Value = "lgbm"
$ws.Range("B8").Value = 0.9035674470457079
$ws.Range("C8").Value = 0.9035674470457079
$ws.Range("D8").Value = 0.90657404350533
$ws.Range("E8").Value = 0.9037620202233
$ws.Range("F8").Value = 0.4533237308263779
$ws.Range("G8").Value = 0.1842628717422485

# Row 9 - mlp (single model row)
$ws.Range("A9").Value = "mlp"
$ws.Range("B9").Value = 0.875139353400223
$ws.Range("C9").Value = 0.875139353400223
$ws.Range("D9").Value = 0.89047978087285
$ws.Range("E9").Value = 0.876280832967023
$ws.Range("F9").Value = 0.3528055310249328
$ws.Range("G9").Value = 0.1454615592956543
